# Generate Report for Handback
# Adds a new handback row (GUID 3b2b63d8-2d6f-4e8e-8162-58619a229075) to the
# "Overview", "zh-cn" and "de-de" worksheets of the handback-status workbook.

function Apply-HyperlinkLook {
    param($rng)
    $f = $rng.Font()
    $f.Name = "Calibri"
    $f.Underline = 2
    $f.Color = 15570276
}

function Apply-DateLook {
    param($rng)
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$wb = $excel.ActiveWorkbook

$guid = "3b2b63d8-2d6f-4e8e-8162-58619a229075"
$handoffHash = "299aaf2c54236d46b12c44e2ac42215ccba29114"

$mdName = "$guid.md"
$zhXlfName = "$guid.$handoffHash.zh-cn.xlf"
$deXlfName = "$guid.$handoffHash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$include = "Include"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e302123b7000bfe42530f72e221829fb876e020a/e2e/$mdName"

$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/835dd33fc86b3584ae22c619172a9e282d63c6e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName"
$zhMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4fc512f52640b349e1d8158494e1cbfaef1f5ee9/e2e/$mdName"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b6d5691b4ffa441542a5012461dca6c89da4288d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName"

$deHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24ab99a96826aaae82378b40c1164abd298f7b20/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName"
$deMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c83c869c1be7f5c902fc9df5cf59bb65ef0e30bf/e2e/$mdName"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a75ff42475ca4e1b2bd12a76bcf10e29da4f70bb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName"

# ---------------------------------------------------------------------------
# Sheet "Overview": add row 4 -> File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
Apply-HyperlinkLook($wsOverview.Range("A4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn": add row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Range("C4").Value = $zhXlfName
$wsZh.Range("D4").Value = "2016-02-17 06:02:06"
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $zhXlfName
$wsZh.Range("G4").Value = "2016-02-17 06:02:49"
$wsZh.Range("H4").Value = $include

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhHandoffUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), $zhMdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), $zhHandbackUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)

Apply-HyperlinkLook($wsZh.Range("A4"))
Apply-HyperlinkLook($wsZh.Range("C4"))
Apply-HyperlinkLook($wsZh.Range("E4"))
Apply-HyperlinkLook($wsZh.Range("F4"))
Apply-DateLook($wsZh.Range("D4"))

# ---------------------------------------------------------------------------
# Sheet "de-de": add row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Range("C4").Value = $deXlfName
$wsDe.Range("D4").Value = "2016-02-17 06:02:17"
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $deXlfName
$wsDe.Range("G4").Value = "2016-02-17 06:03:06"
$wsDe.Range("H4").Value = $include

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deHandoffUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), $deMdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), $deHandbackUrl, [Type]::Missing, [Type]::Missing, $deXlfName)

Apply-HyperlinkLook($wsDe.Range("A4"))
Apply-HyperlinkLook($wsDe.Range("C4"))
Apply-HyperlinkLook($wsDe.Range("E4"))
Apply-HyperlinkLook($wsDe.Range("F4"))
Apply-DateLook($wsDe.Range("D4"))
